$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.558.20'
$ws.Range("E2").Value = '  +1.57%  '

$ws.Range("D3").Value = '1.473.22'
$ws.Range("E3").Value = '  +2.24%  '

$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").Value = '''0.9584'
$ws.Range("E5").Value = '  +3.58%  '

$ws.Range("D6").Value = '''276.92'
$ws.Range("E6").Value = '  +1.27%  '

$ws.Range("D7").Value = '''0.3520'
$ws.Range("E7").Value = '  -3.44%  '

$ws.Range("D8").Value = '''0.3069'
$ws.Range("E8").Value = '  +0.19%  '

$ws.Range("D9").Value = '''1.081'
$ws.Range("E9").Value = '  +6.20%  '

$ws.Range("D10").Value = '''39.31'
$ws.Range("E10").Value = '  -0.31%  '

$ws.Range("D11").Value = '''0.06629'
$ws.Range("E11").Value = '  +2.00%  '

$ws.Range("E12").Value = '  +0.57%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''5.448'
$ws.Range("E13").Value = '  +2.02%  '

$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").Value = '''18.02'
$ws.Range("E14").Value = '  +3.36%  '

$ws.Range("D15").Value = '''6.162'
$ws.Range("E15").Value = '  +1.97%  '

$ws.Range("D16").Value = '''0.9573'
$ws.Range("E16").Value = '  +1.44%  '

$ws.Range("D17").Value = '''0.00001013'
$ws.Range("E17").Value = '  +0.41%  '

$ws.Range("D18").Value = '1.472.24'
$ws.Range("E18").Value = '  +2.34%  '

$ws.Range("D19").Value = '''0.05962'
$ws.Range("E19").Value = '  +5.48%  '

$ws.Range("D20").Value = '''68.69'
$ws.Range("E20").Value = '  +0.45%  '

$ws.Range("D21").Value = '''5.472'
$ws.Range("E21").Value = '  +2.15%  '

$ws.Range("D22").Value = '''14.44'
$ws.Range("E22").Value = '  +1.60%  '

$ws.Range("D23").Value = '''11.13'
$ws.Range("E23").Value = '  +3.34%  '

$ws.Range("D24").Value = '''2.280'
$ws.Range("E24").Value = '  +1.59%  '

$ws.Range("D25").Value = '20.578.22'
$ws.Range("E25").Value = '  +1.57%  '

$ws.Range("D26").Value = '''146.25'
$ws.Range("E26").Value = '  +4.23%  '

$ws.Range("D27").Value = '''2.073'
$ws.Range("E27").Value = '  +2.36%  '

$ws.Range("D28").Value = '''17.11'
$ws.Range("E28").Value = '  +1.40%  '

$ws.Range("D29").Value = '1.635.41'
$ws.Range("E29").Value = '  +2.81%  '

$ws.Range("D30").Value = '''114.19'
$ws.Range("E30").Value = '  +3.66%  '

$ws.Range("D31").Value = '''3.934'
$ws.Range("E31").Value = '  -2.08%  '

$ws.Range("D32").Value = '''4.919'
$ws.Range("E32").Value = '  +2.81%  '

$ws.Range("D33").Value = '''0.07887'
$ws.Range("E33").Value = '  +2.80%  '

$ws.Range("D34").Value = '''0.7920'
$ws.Range("E34").Value = '  +2.29%  '

$ws.Range("D35").Value = '''1.198'
$ws.Range("E35").Value = '  +7.74%  '

$ws.Range("D36").Value = '''1.435'
$ws.Range("E36").Value = '  -1.10%  '

$ws.Range("D37").Value = '''0.05668'
$ws.Range("E37").Value = '  +0.34%  '

$ws.Range("D38").Value = '''4.680'
$ws.Range("E38").Value = '  +1.11%  '

$ws.Range("D39").Value = '''0.9584'
$ws.Range("E39").Value = '  +2.37%  '

$ws.Range("D40").Value = '''0.02014'
$ws.Range("E40").Value = '  +1.40%  '

$ws.Range("D41").Value = '''10.21'
$ws.Range("E41").Value = '  +0.81%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '''0.1841'
$ws.Range("E42").Value = '  +0.50%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''7.346'
$ws.Range("E43").Value = '  +5.62%  '

$ws.Range("B44").Value = 'PancakeSwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D44").Value = '''3.509'
$ws.Range("E44").Value = '  +1.17%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '''0.5203'
$ws.Range("E45").Value = '  +0.54%  '

$ws.Range("D46").Value = '''12.02'
$ws.Range("E46").Value = '  +3.12%  '

$ws.Range("D47").Value = '''120.00'
$ws.Range("E47").Value = '  +5.08%  '

$ws.Range("D48").Value = '''0.5150'
$ws.Range("E48").Value = '  +1.84%  '

$ws.Range("D49").Value = '''1.804'
$ws.Range("E49").Value = '  +4.83%  '

$ws.Range("D50").Value = '''0.06400'
$ws.Range("E50").Value = '  +0.67%  '

$ws.Range("D51").Value = '''0.9929'
$ws.Range("E51").Value = '  +0.43%  '
